# Formed the consolidated report:
# Compute the "Absent" column (H) for each attendance row as the complement
# of the "Real" column (E): a genuine (real) attendance means the student
# was present (H = 0); anything else (no record, duplicate-only, invalid
# swipe, etc.) means the student counts as absent for that date (H = 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows start at row 3 (row 1 = header, row 2 = roll/name info) and
# run through row 21, matching the sheet's used range (A1:H21).
$firstRow = 3
$lastRow = 21

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $realAttendance = $ws.Cells.Item($r, 5).Value2   # column E - "Real"
    if ($realAttendance -eq 1) {
        $ws.Cells.Item($r, 8).Value = 0              # column H - present
    } else {
        $ws.Cells.Item($r, 8).Value = 1              # column H - absent
    }
}
